$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "21.684.50"
$ws.Range("E2").NumberFormat = "@"
$ws.Range("E2").Value = "  -1.50%  "

$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "1.534.78"
$ws.Range("E3").NumberFormat = "@"
$ws.Range("E3").Value = "  -1.26%  "

$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = "1.001"
$ws.Range("E4").NumberFormat = "@"
$ws.Range("E4").Value = "  +0.12%  "

$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "289.18"
$ws.Range("E6").NumberFormat = "@"
$ws.Range("E6").Value = "  +0.87%  "

$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "0.3934"
$ws.Range("E7").NumberFormat = "@"
$ws.Range("E7").Value = "  +3.88%  "

$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.3153"
$ws.Range("E8").NumberFormat = "@"
$ws.Range("E8").Value = "  -2.59%  "

$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "42.40"
$ws.Range("E9").NumberFormat = "@"
$ws.Range("E9").Value = "  +2.93%  "

$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "0.07165"
$ws.Range("E10").NumberFormat = "@"
$ws.Range("E10").Value = "  -1.84%  "

$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "1.046"
$ws.Range("E11").NumberFormat = "@"
$ws.Range("E11").Value = "  -6.92%  "

$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "1.001"
$ws.Range("E12").NumberFormat = "@"
$ws.Range("E12").Value = "  +0.13%  "

$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "5.604"
$ws.Range("E13").NumberFormat = "@"
$ws.Range("E13").Value = "  -1.99%  "

$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "18.43"
$ws.Range("E14").NumberFormat = "@"
$ws.Range("E14").Value = "  -4.63%  "

$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "6.591"
$ws.Range("E15").NumberFormat = "@"
$ws.Range("E15").Value = "  -2.78%  "

$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "1.535.02"
$ws.Range("E16").NumberFormat = "@"
$ws.Range("E16").Value = "  -1.36%  "

$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "0.00001097"
$ws.Range("E17").NumberFormat = "@"
$ws.Range("E17").Value = "  +1.04%  "

$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "0.06579"
$ws.Range("E18").NumberFormat = "@"
$ws.Range("E18").Value = "  -0.57%  "

$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "82.89"
$ws.Range("E19").NumberFormat = "@"
$ws.Range("E19").Value = "  -2.46%  "

$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "1.000"
$ws.Range("E20").NumberFormat = "@"
$ws.Range("E20").Value = "  +0.20%  "

$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "6.112"
$ws.Range("E21").NumberFormat = "@"
$ws.Range("E21").Value = "  -4.68%  "

$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "15.34"
$ws.Range("E22").NumberFormat = "@"
$ws.Range("E22").Value = "  -3.78%  "

$ws.Range("E23").NumberFormat = "@"
$ws.Range("E23").Value = "  -5.39%  "

$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "2.382"
$ws.Range("E24").NumberFormat = "@"
$ws.Range("E24").Value = "  +4.24%  "

$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "21.683.58"
$ws.Range("E25").NumberFormat = "@"
$ws.Range("E25").Value = "  -1.55%  "

$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "2.342"
$ws.Range("E26").NumberFormat = "@"
$ws.Range("E26").Value = "  -7.54%  "

$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "146.44"

$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "18.29"
$ws.Range("E28").NumberFormat = "@"
$ws.Range("E28").Value = "  -2.65%  "

$ws.Range("E29").NumberFormat = "@"
$ws.Range("E29").Value = "  -0.44%  "

$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "1.708.66"
$ws.Range("E30").NumberFormat = "@"
$ws.Range("E30").Value = "  -1.28%  "

$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "116.84"
$ws.Range("E31").NumberFormat = "@"
$ws.Range("E31").Value = "  -2.61%  "

$ws.Range("B32").Value = "ImmutableX"
$ws.Range("C32").Value = "https://coinranking.com/coin/Z96jIvLU7+immutablex-imx"
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "0.9594"
$ws.Range("E32").NumberFormat = "@"
$ws.Range("E32").Value = "  -13.29%  "

$ws.Range("B33").Value = "Filecoin"
$ws.Range("C33").Value = "https://coinranking.com/coin/ymQub4fuB+filecoin-fil"
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "5.853"
$ws.Range("E33").NumberFormat = "@"
$ws.Range("E33").Value = "  -1.22%  "

$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "0.08159"
$ws.Range("E34").NumberFormat = "@"
$ws.Range("E34").Value = "  +0.20%  "

$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "8.646"
$ws.Range("E35").NumberFormat = "@"
$ws.Range("E35").Value = "  -6.22%  "

$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "0.06053"
$ws.Range("E36").NumberFormat = "@"
$ws.Range("E36").Value = "  -2.03%  "

$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "5.090"
$ws.Range("E37").NumberFormat = "@"
$ws.Range("E37").Value = "  -2.76%  "

$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "0.02193"

$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "0.2019"
$ws.Range("E39").NumberFormat = "@"
$ws.Range("E39").Value = "  -4.45%  "

$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "1.436"
$ws.Range("E40").NumberFormat = "@"
$ws.Range("E40").Value = "  -13.01%  "

$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "1.180"
$ws.Range("E41").NumberFormat = "@"
$ws.Range("E41").Value = "  -3.10%  "

$ws.Range("E42").NumberFormat = "@"
$ws.Range("E42").Value = "  +0.12%  "

$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "10.66"
$ws.Range("E43").NumberFormat = "@"
$ws.Range("E43").Value = "  -1.99%  "

$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "0.5706"
$ws.Range("E44").NumberFormat = "@"
$ws.Range("E44").Value = "  -3.80%  "

$ws.Range("B45").Value = "EnergySwap"
$ws.Range("C45").Value = "https://coinranking.com/coin/SbWqqTui-+energyswap-ens"
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "13.03"
$ws.Range("E45").NumberFormat = "@"
$ws.Range("E45").Value = "  -3.67%  "

$ws.Range("B46").Value = "PancakeSwap"
$ws.Range("C46").Value = "https://coinranking.com/coin/ncYFcP709+pancakeswap-cake"
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "3.724"
$ws.Range("E46").NumberFormat = "@"
$ws.Range("E46").Value = "  +0.18%  "

$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "0.5464"
$ws.Range("E47").NumberFormat = "@"
$ws.Range("E47").Value = "  -4.64%  "

$ws.Range("E48").NumberFormat = "@"
$ws.Range("E48").Value = "  +0.63%  "

$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "116.00"
$ws.Range("E49").NumberFormat = "@"
$ws.Range("E49").Value = "  -3.03%  "

$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "1.856"
$ws.Range("E50").NumberFormat = "@"
$ws.Range("E50").Value = "  -3.98%  "

$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "0.06684"
$ws.Range("E51").NumberFormat = "@"
$ws.Range("E51").Value = "  -2.92%  "
